$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Fabrica 3 (Lego)
$ws.Range("A4").Value = "Fabrica 3:"
$ws.Range("B4").Value = "Lego"
$ws.Range("C4").Value = "EE.UU."
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 10

# Row 5: Fabrica 4 (Bandai Namco)
$ws.Range("A5").Value = "Fabrica 4:"
$ws.Range("B5").Value = "Bandai Namco"
$ws.Range("C5").Value = "China"
$ws.Range("D5").Value = 9000
$ws.Range("E5").Value = 8

# Match the label styling already used on A2/A3 (bold, bordered, centered/top)
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A4").Value = "Fabrica 3:"
$ws.Range("A5").Value = "Fabrica 4:"
$excel.CutCopyMode = $false
